$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 7011466872
$ws.Range("B3").Value = "SinUsername"
$ws.Range("C3").Value = "2025-09-09 21:51:24"
